# Atualizei dados da ADD
# - Corrige o faturamento de 11/06/2025 e 12/06/2025
# - Insere um novo registro para o dia 13/06/2025 (valor 6825.70)
#   empurrando os registros seguintes uma linha para baixo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrige os valores de total_venda existentes para os dias 11 e 12/06/2025
$ws.Range("B9").Value  = 20439.56
$ws.Range("B10").Value = 17899.38

# Insere uma nova linha em branco na posição 11 (desloca as linhas 11..70 para 12..71)
$ws.Rows(11).Insert()

# Preenche a nova linha com o registro do dia 13/06/2025
$ws.Range("A11").Value = 13
$ws.Range("B11").Value = 6825.7
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = 2025
$ws.Range("E11").Value = "06/2025"
